# #5: fund, bonds, otherbonds, antique done
# Normalize/complete the "具有相當價值之財產" (property of considerable value)
# sheet (the 6th worksheet) by turning row 1 into a real header row and by
# appending the standard trailing metadata columns (property_category,
# category, date, legislator_name, legislator_id, source_file, index) to
# every data row. Also fixes two garbled quantity values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# ---- Row 1: turn the stray duplicated data row into real column headers ----
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "quantity"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "total"
$ws.Cells.Item(1,6).Value = "property_category"
$ws.Cells.Item(1,7).Value = "category"
$ws.Cells.Item(1,8).Value = "date"
$ws.Cells.Item(1,9).Value = "legislator_name"
$ws.Cells.Item(1,10).Value = "legislator_id"
$ws.Cells.Item(1,11).Value = "source_file"
$ws.Cells.Item(1,12).Value = "index"

# ---- Fix two mis-OCR'd / corrupted quantity values ----
$ws.Cells.Item(5,3).Value = "10件"
$ws.Cells.Item(7,3).Value = "10件"

# ---- Rows 2-10: append the standard metadata columns ----
$indices = @(88, 89, 90, 91, 92, 93, 94, 95, 96)
for ($i = 0; $i -lt $indices.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r,6).Value  = "otherbonds"
    $ws.Cells.Item($r,7).Value  = "normal"
    $ws.Cells.Item($r,8).Value  = "2012-04-27"
    $ws.Cells.Item($r,9).Value  = "李貴敏"
    $ws.Cells.Item($r,10).Value = 1739
    $ws.Cells.Item($r,11).Value = "tmp59331"
    $ws.Cells.Item($r,12).Value = $indices[$i]
}
